# Update the "dSF" column (F) values for specific rows on Sheet1
# as part of a re-pull/push of data and mean calculation refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = -3
    11 = 2
    19 = 2
    20 = 1
    24 = 2
    27 = -2
    33 = 2
    40 = 2
    41 = 2
    44 = -2
    49 = -1
    51 = 3
    70 = -4
    71 = -3
    74 = -5
    78 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
